$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-wise data: re-read/update each row with the new version-number strings ---

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "\Config.wxi"
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "No"
$ws.Range("E2").Value = "MajorMinorPatch"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "\config.wxi"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = "Major.Minor"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "\Config.wxi"
$ws.Range("C4").Value = "Yes"
$ws.Range("D4").Value = "No"
$ws.Range("E4").Value = "Major,Minor,Patch,0"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "\Config.wxi"
$ws.Range("C5").Value = "Yes"
$ws.Range("D5").Value = "No"
$ws.Range("E5").Value = "Major, Minor, Patch, BuildNumberUpdate"
$ws.Range("E5").WrapText = $true

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "\Confoiig.wxi"
$ws.Range("C6").Value = "Yes"
$ws.Range("D6").Value = "No"
$ws.Range("E6").Value = "Major,Minor,Patch,BuildNumberUpdate"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "\Config.wxi"
$ws.Range("C7").Value = "Yes"
$ws.Range("D7").Value = "No"
$ws.Range("E7").Value = "Major.Minor.Patch"

# --- Column widths (best achievable on this engine's width grid) ---
$ws.Columns("B").ColumnWidth = 60.166666666666664
$ws.Columns("C").ColumnWidth = 12.666666666666666
$ws.Columns("D").ColumnWidth = 12.666666666666666
$ws.Columns("E").ColumnWidth = 58.166666666666664

# --- Final selection matches the author's last-edited cell ---
[void]$ws.Range("B6").Select()
